$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.206.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +6.24%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.919.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.85%  "

$ws.Range("E4").Value = "  -0.60%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.31%  "

$ws.Range("E6").Value = "  -0.61%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5230"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.28%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4093"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.90%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08522"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.50%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.04"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.47%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.130"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.45%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.10"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +13.94%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.460"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.73%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.892.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.54%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.408"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.54%  "

$ws.Range("E16").Value = "  -0.65%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "95.37"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.77%  "

$ws.Range("E18").Value = "  +1.65%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06702"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.33%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9993"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.040"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.229.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.23%  "

$ws.Range("E24").Value = "  +3.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.221"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.26%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.124.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.28%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.44%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.42%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.421"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.64%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.38%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.103"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.81%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1073"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.91%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.029"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.58%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.604"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.46%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02497"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.49%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06602"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.78%  "

$ws.Range("E37").Value = "  +3.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.237"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.88%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.186"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.54%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.47%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.825"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.29%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6554"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.45%  "

$ws.Range("E43").Value = "  +0.80%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6178"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.64%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.752"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.17%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.096"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.23%  "

$ws.Range("E48").Value = "  +3.39%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.22%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.164"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.01%  "

$ws.Range("E51").Value = "  +5.29%  "

